$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.853.80"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +4.91%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.365.03"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +2.82%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.20%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'548.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +2.59%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'133.34"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +1.97%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  +1.67%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'2.362.85"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +2.76%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  +2.18%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  +1.71%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'  +1.08%  "
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'  +1.82%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'24.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +2.92%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'2.788.43"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +3.86%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'60.722.10"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +4.76%  "
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = "'  +2.14%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.359.46"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +2.95%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'10.76"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +2.32%  "
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'  +9.09%  "
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'  -0.51%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'316.56"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.97%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +0.08%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'63.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +1.37%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  +3.51%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  +0.31%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +0.37%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'1.36"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +5.49%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  +3.07%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'172.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +0.93%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  +2.49%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +10.48%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'5.93"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +2.74%  "
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = "'  +16.61%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'0.384"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +1.74%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'18.10"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +1.86%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  +0.02%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -0.14%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'4.19"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +7.65%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'318.39"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +10.40%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  +3.70%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'38.25"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +0.53%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'143.86"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +2.05%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  +2.23%  "
$ws.Range("E44").ClearFormats()
$ws.Range("E45").Value = "'  +1.08%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'19.41"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +7.47%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'0.0501"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +0.98%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.565"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +1.79%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'  +2.03%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0₆0215"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +6.98%  "
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'  +1.09%  "
$ws.Range("E51").ClearFormats()
